$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Date_Time formulas in column C (rows 2-30) to match the new
# `from_arbin_to_datetime` function: a running cumulative sum seeded from a
# fixed epoch-like constant instead of the old day-fraction calculation.
$ws.Range("C2").Formula = "=16566949984405000+10000000*E2"
$ws.Range("C3").Formula = "=C2+10000000*E3"
$ws.Range("C4:C30").Formula = "=C3+10000000*E4"

# Column C (Date_Time) is widened so the new, longer formula/value is fully visible.
$ws.Columns.Item(3).ColumnWidth = 62.5

# Selection narrowed from the whole fill range down to the first formula cell.
[void]$ws.Range("C3").Select()
